$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 777.0476
$ws.Range("I53").Value = 289.1111
$ws.Range("J53").Value = 1143
$ws.Range("K53").Value = 289.1111
$ws.Range("L53").Value = 1143
$ws.Range("M53").Value = 347.8889
$ws.Range("N53").Value = -2417
$ws.Range("H55").Value = 339.07144
$ws.Range("I55").Value = 337.33334
$ws.Range("J55").Value = 349.5
$ws.Range("K55").Value = 337.33334
$ws.Range("L55").Value = 349.5
$ws.Range("M55").Value = -123.33334
$ws.Range("N55").Value = -777.5
$ws.Range("H98").Value = 2305.0588
$ws.Range("I98").Value = 2858.7144
$ws.Range("J98").Value = 1917.5
$ws.Range("K98").Value = 2858.7144
$ws.Range("L98").Value = 1917.5
$ws.Range("M98").Value = -1360.7144
$ws.Range("N98").Value = -4913.5
$ws.Range("H122").Value = 2305.0588
$ws.Range("I122").Value = 2858.7144
$ws.Range("J122").Value = 1917.5
$ws.Range("K122").Value = 8576.143199999999
$ws.Range("L122").Value = 5752.5
$ws.Range("M122").Value = -6126.143199999999
$ws.Range("N122").Value = -10652.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H130").Value = 200041400
$ws.Range("J130").Value = 200041400
$ws.Range("L130").Value = 200041400
$ws.Range("N130").Value = -200051440
$ws.Range("H132").Value = 4832.821
$ws.Range("I132").Value = 4014.7446
$ws.Range("J132").Value = 6755.3
$ws.Range("K132").Value = 12044.2338
$ws.Range("L132").Value = 20265.9
$ws.Range("M132").Value = -9514.2338
$ws.Range("N132").Value = -25325.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 12518
$ws.Range("I28").Value = 10517.75
$ws.Range("J28").Value = 20519
$ws.Range("K28").Value = 10517.75
$ws.Range("L28").Value = 20519
$ws.Range("M28").Value = -10325.75
$ws.Range("N28").Value = -20903
$ws.Range("H45").Value = 4365
$ws.Range("I45").Value = 2235.3333
$ws.Range("J45").Value = 8624.333000000001
$ws.Range("K45").Value = 2235.3333
$ws.Range("L45").Value = 8624.333000000001
$ws.Range("M45").Value = -1858.3333
$ws.Range("N45").Value = -9378.333000000001
$ws.Range("H61").Value = 3377.5
$ws.Range("I61").Value = 3402
$ws.Range("J61").Value = 3369.3333
$ws.Range("K61").Value = 3402
$ws.Range("L61").Value = 3369.3333
$ws.Range("M61").Value = -3190
$ws.Range("N61").Value = -3793.3333
$ws.Range("H99").Value = 12518
$ws.Range("I99").Value = 10517.75
$ws.Range("J99").Value = 20519
$ws.Range("K99").Value = 10517.75
$ws.Range("L99").Value = 20519
$ws.Range("M99").Value = -7522.75
$ws.Range("N99").Value = -26509
$ws.Range("H133").Value = 27736
$ws.Range("I133").Value = 24998
$ws.Range("J133").Value = 28127.143
$ws.Range("K133").Value = 24998
$ws.Range("L133").Value = 28127.143
$ws.Range("M133").Value = -22468
$ws.Range("N133").Value = -33187.143
$ws.Range("H136").Value = 3377.5
$ws.Range("I136").Value = 3402
$ws.Range("J136").Value = 3369.3333
$ws.Range("K136").Value = 10206
$ws.Range("L136").Value = 10107.9999
$ws.Range("M136").Value = -7656
$ws.Range("N136").Value = -15207.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2205.7144
$ws.Range("I105").Value = 1838.8889
$ws.Range("J105").Value = 2866
$ws.Range("K105").Value = 1838.8889
$ws.Range("L105").Value = 2866
$ws.Range("M105").Value = -91.88889999999992
$ws.Range("N105").Value = -6360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4426.375
$ws.Range("I16").Value = 4568.5
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 4568.5
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -4281.5
$ws.Range("N16").Value = -4574
$ws.Range("H31").Value = 5466177.5
$ws.Range("I31").Value = 1224.0667
$ws.Range("K31").Value = 1224.0667
$ws.Range("M31").Value = -929.0667000000001
$ws.Range("H34").Value = 5466177.5
$ws.Range("I34").Value = 1224.0667
$ws.Range("K34").Value = 1224.0667
$ws.Range("M34").Value = -1022.0667
$ws.Range("H58").Value = 964800.4
$ws.Range("I58").Value = 1733.6072
$ws.Range("J58").Value = 2088378.2
$ws.Range("K58").Value = 1733.6072
$ws.Range("L58").Value = 2088378.2
$ws.Range("M58").Value = -1530.6072
$ws.Range("N58").Value = -2088784.2
$ws.Range("H113").Value = 4426.375
$ws.Range("I113").Value = 4568.5
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 4568.5
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -2398.5
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 55557076
$ws.Range("I122").Value = 83334264
$ws.Range("J122").Value = 2699.8333
$ws.Range("K122").Value = 250002792
$ws.Range("L122").Value = 8099.499899999999
$ws.Range("M122").Value = -250000342
$ws.Range("N122").Value = -12999.4999
$ws.Range("H136").Value = 964800.4
$ws.Range("I136").Value = 1733.6072
$ws.Range("J136").Value = 2088378.2
$ws.Range("K136").Value = 5200.821599999999
$ws.Range("L136").Value = 6265134.6
$ws.Range("M136").Value = -2650.821599999999
$ws.Range("N136").Value = -6270234.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 129.6875
$ws.Range("I38").Value = 98.75
$ws.Range("J38").Value = 160.625
$ws.Range("K38").Value = 296.25
$ws.Range("L38").Value = 481.875
$ws.Range("M38").Value = 50.75
$ws.Range("N38").Value = -1175.875
$ws.Range("H125").Value = 5714.3335
$ws.Range("I125").Value = 2628.75
$ws.Range("J125").Value = 6836.364
$ws.Range("K125").Value = 7886.25
$ws.Range("L125").Value = 20509.092
$ws.Range("M125").Value = -2966.25
$ws.Range("N125").Value = -30349.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5498151.5
$ws.Range("I102").Value = 8931493
$ws.Range("J102").Value = 4805.6
$ws.Range("K102").Value = 8931493
$ws.Range("L102").Value = 4805.6
$ws.Range("M102").Value = -8929871
$ws.Range("N102").Value = -8049.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2298.318
$ws.Range("I7").Value = 2150.5
$ws.Range("K7").Value = 2150.5
$ws.Range("M7").Value = -2038.5
$ws.Range("H22").Value = 1072.4
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 1169.3334
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 1169.3334
$ws.Range("M22").Value = 95
$ws.Range("N22").Value = -1759.3334
$ws.Range("H27").Value = 1072.4
$ws.Range("I27").Value = 200
$ws.Range("J27").Value = 1169.3334
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 1169.3334
$ws.Range("M27").Value = -93
$ws.Range("N27").Value = -1383.3334
$ws.Range("H126").Value = 2298.318
$ws.Range("I126").Value = 2150.5
$ws.Range("K126").Value = 6451.5
$ws.Range("M126").Value = -3981.5
$ws.Range("H132").Value = 27056264
$ws.Range("I132").Value = 37075824
$ws.Range("K132").Value = 111227472
$ws.Range("M132").Value = -111224942
$ws.Range("H136").Value = 33336140
$ws.Range("I136").Value = 45455920
$ws.Range("J136").Value = 6750
$ws.Range("K136").Value = 136367760
$ws.Range("L136").Value = 20250
$ws.Range("M136").Value = -136365210
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 79804
$ws.Range("I15").Value = 86668.664
$ws.Range("J15").Value = 69507
$ws.Range("K15").Value = 86668.664
$ws.Range("L15").Value = 69507
$ws.Range("M15").Value = -86380.664
$ws.Range("N15").Value = -70083
